# Actualización automática 2025-08-15 16:00:08
# Insert a new client row ("FERRETERIAS FERRIGONZ SA") above the
# "JARAMILLO CARVAJAL NICOLAS ESTEBAN" row (row 10) on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing every later
# client row down by one, and refresh the "X de N" compliance-count
# footer on "VENTAS POR GRUPO" to reflect the new total of 21 clients.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(10).Insert()

$ws1.Range("A10").Value = "HIDALGO HIDALGO PEDRO GUSTAVO"
$ws1.Range("B10").Value = "FERRETERIAS FERRIGONZ SA"
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $ws1.Range($col + "10").Value = 0
}

# Footer row (previously row 22, now row 23): bump the "X de 20"
# counters to "X de 21" now that there are 21 clients.
$ws1.Range("C23").Value = "1 de 21"
$ws1.Range("D23").Value = "0 de 21"
$ws1.Range("E23").Value = "0 de 21"
$ws1.Range("F23").Value = "0 de 21"
$ws1.Range("G23").Value = "0 de 21"
$ws1.Range("H23").Value = "1 de 21"
$ws1.Range("I23").Value = "1 de 21"
$ws1.Range("J23").Value = "0 de 21"
$ws1.Range("K23").Value = "0 de 21"
$ws1.Range("L23").Value = "0 de 21"
$ws1.Range("M23").Value = "2 de 21"
$ws1.Range("N23").Value = "0 de 21"
$ws1.Range("O23").Value = "0 de 21"
$ws1.Range("P23").Value = "0 de 21"
$ws1.Range("Q23").Value = "0 de 21"
$ws1.Range("R23").Value = "0 de 21"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(10).Insert()

$ws2.Range("A10").Value = "HIDALGO HIDALGO PEDRO GUSTAVO"
$ws2.Range("B10").Value = "FERRETERIAS FERRIGONZ SA"
$cols2 = @("C","D","E","F","G")
foreach ($col in $cols2) {
    $ws2.Range($col + "10").Value = 0
}
